$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 714352.9
$ws.Range("I9").Value = 1000054
$ws.Range("K9").Value = 1000054
$ws.Range("M9").Value = -999885
$ws.Range("H38").Value = 3093.4546
$ws.Range("I38").Value = 868.5
$ws.Range("J38").Value = 6987.125
$ws.Range("K38").Value = 2605.5
$ws.Range("L38").Value = 20961.375
$ws.Range("M38").Value = -2233.5
$ws.Range("N38").Value = -21705.375
$ws.Range("H58").Value = 19231810
$ws.Range("I58").Value = 22727594
$ws.Range("K58").Value = 68182782
$ws.Range("M58").Value = -68182632
$ws.Range("H64").Value = 8482.725
$ws.Range("J64").Value = 9565.174000000001
$ws.Range("L64").Value = 9565.174000000001
$ws.Range("N64").Value = -10061.174
$ws.Range("H67").Value = 8482.725
$ws.Range("J67").Value = 9565.174000000001
$ws.Range("L67").Value = 9565.174000000001
$ws.Range("N67").Value = -11281.174
$ws.Range("H98").Value = 16895.2
$ws.Range("I98").Value = 27495.834
$ws.Range("J98").Value = 994.25
$ws.Range("K98").Value = 27495.834
$ws.Range("L98").Value = 994.25
$ws.Range("M98").Value = -25997.834
$ws.Range("N98").Value = -3990.25
$ws.Range("H100").Value = 4461.778
$ws.Range("I100").Value = 1431.3
$ws.Range("J100").Value = 8249.875
$ws.Range("K100").Value = 1431.3
$ws.Range("L100").Value = 8249.875
$ws.Range("M100").Value = -890.3
$ws.Range("N100").Value = -9331.875
$ws.Range("H113").Value = 4857.5713
$ws.Range("I113").Value = 2799.4
$ws.Range("K113").Value = 2799.4
$ws.Range("M113").Value = 454.5999999999999
$ws.Range("H122").Value = 16895.2
$ws.Range("I122").Value = 27495.834
$ws.Range("J122").Value = 994.25
$ws.Range("K122").Value = 82487.50199999999
$ws.Range("L122").Value = 2982.75
$ws.Range("M122").Value = -80037.50199999999
$ws.Range("N122").Value = -7882.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2659.451
$ws.Range("I32").Value = 2387.1914
$ws.Range("K32").Value = 2387.1914
$ws.Range("M32").Value = -2100.1914
$ws.Range("H45").Value = 2155.5715
$ws.Range("I45").Value = 2238
$ws.Range("J45").Value = 1949.5
$ws.Range("K45").Value = 2238
$ws.Range("L45").Value = 1949.5
$ws.Range("M45").Value = -1861
$ws.Range("N45").Value = -2703.5
$ws.Range("H97").Value = 372.3913
$ws.Range("I97").Value = 192.5
$ws.Range("K97").Value = 192.5
$ws.Range("M97").Value = 303.5
$ws.Range("H132").Value = 3757.375
$ws.Range("I132").Value = 3757.375
$ws.Range("K132").Value = 11272.125
$ws.Range("M132").Value = -8742.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1889.2
$ws.Range("I20").Value = 1332.3334
$ws.Range("K20").Value = 1332.3334
$ws.Range("M20").Value = -1085.3334
$ws.Range("H86").Value = 705239.9
$ws.Range("I86").Value = 1449.9474
$ws.Range("J86").Value = 2376741
$ws.Range("K86").Value = 1449.9474
$ws.Range("L86").Value = 2376741
$ws.Range("M86").Value = -326.9474
$ws.Range("N86").Value = -2378987
$ws.Range("H89").Value = 705239.9
$ws.Range("I89").Value = 1449.9474
$ws.Range("J89").Value = 2376741
$ws.Range("K89").Value = 7249.737
$ws.Range("L89").Value = 11883705
$ws.Range("M89").Value = -1633.737
$ws.Range("N89").Value = -11894937
$ws.Range("H92").Value = 33999.4
$ws.Range("J92").Value = 33999.4
$ws.Range("L92").Value = 33999.4
$ws.Range("N92").Value = -38991.4
$ws.Range("H134").Value = 2123.44
$ws.Range("I134").Value = 2094.818
$ws.Range("K134").Value = 6284.454000000001
$ws.Range("M134").Value = -3749.454000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1728.3
$ws.Range("I122").Value = 1923
$ws.Range("J122").Value = 949.5
$ws.Range("K122").Value = 5769
$ws.Range("L122").Value = 2848.5
$ws.Range("M122").Value = -3319
$ws.Range("N122").Value = -7748.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4303322
$ws.Range("I4").Value = 976126.8
$ws.Range("K4").Value = 2928380.4
$ws.Range("M4").Value = -2928268.4
$ws.Range("H5").Value = 1114.909
$ws.Range("I5").Value = 744
$ws.Range("K5").Value = 2232
$ws.Range("M5").Value = -2120
$ws.Range("H58").Value = 465
$ws.Range("J58").Value = 750
$ws.Range("L58").Value = 2250
$ws.Range("N58").Value = -2506
$ws.Range("H68").Value = 1005.82355
$ws.Range("I68").Value = 400
$ws.Range("J68").Value = 1043.6875
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 3131.0625
$ws.Range("M68").Value = -389
$ws.Range("N68").Value = -4753.0625
$ws.Range("H71").Value = 1005.82355
$ws.Range("I71").Value = 400
$ws.Range("J71").Value = 1043.6875
$ws.Range("K71").Value = 3600
$ws.Range("L71").Value = 9393.1875
$ws.Range("M71").Value = 456
$ws.Range("N71").Value = -17505.1875
$ws.Range("H107").Value = 1238.4117
$ws.Range("J107").Value = 855.7273
$ws.Range("L107").Value = 2567.1819
$ws.Range("N107").Value = -6407.1819
$ws.Range("H135").Value = 1114.909
$ws.Range("I135").Value = 744
$ws.Range("K135").Value = 6696
$ws.Range("M135").Value = -4161
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 11969
$ws.Range("I68").Value = 11969
$ws.Range("K68").Value = 11969
$ws.Range("M68").Value = -11158
$ws.Range("H71").Value = 11969
$ws.Range("I71").Value = 11969
$ws.Range("K71").Value = 35907
$ws.Range("M71").Value = -31851
$ws.Range("H97").Value = 404
$ws.Range("I97").Value = 489.875
$ws.Range("J97").Value = 232.25
$ws.Range("K97").Value = 489.875
$ws.Range("L97").Value = 232.25
$ws.Range("M97").Value = 6.125
$ws.Range("N97").Value = -1224.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1802.9565
$ws.Range("I22").Value = 1984.0625
$ws.Range("J22").Value = 1389
$ws.Range("K22").Value = 1984.0625
$ws.Range("L22").Value = 1389
$ws.Range("M22").Value = -1689.0625
$ws.Range("N22").Value = -1979
$ws.Range("H27").Value = 1802.9565
$ws.Range("I27").Value = 1984.0625
$ws.Range("J27").Value = 1389
$ws.Range("K27").Value = 1984.0625
$ws.Range("L27").Value = 1389
$ws.Range("M27").Value = -1877.0625
$ws.Range("N27").Value = -1603
$ws.Range("H69").Value = 199945
$ws.Range("J69").Value = 199945
$ws.Range("L69").Value = 199945
$ws.Range("N69").Value = -201567
$ws.Range("H72").Value = 199945
$ws.Range("J72").Value = 199945
$ws.Range("L72").Value = 599835
$ws.Range("N72").Value = -607947
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H82").Value = 1998.8
$ws.Range("I82").Value = 841.8125
$ws.Range("K82").Value = 841.8125
$ws.Range("M82").Value = -480.8125
$ws.Range("H85").Value = 1998.8
$ws.Range("I85").Value = 841.8125
$ws.Range("K85").Value = 841.8125
$ws.Range("M85").Value = 406.1875
$ws.Range("H93").Value = 3625.9355
$ws.Range("J93").Value = 5163.6
$ws.Range("L93").Value = 5163.6
$ws.Range("N93").Value = -7659.6
$ws.Range("H122").Value = 6997.0557
$ws.Range("I122").Value = 7288
$ws.Range("K122").Value = 21864
$ws.Range("M122").Value = -19414
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 50000612
$ws.Range("I5").Value = 1222
$ws.Range("J5").Value = 100000000
$ws.Range("K5").Value = 1222
$ws.Range("L5").Value = 100000000
$ws.Range("M5").Value = -1110
$ws.Range("N5").Value = -100000224
$ws.Range("H63").Value = 19371.25
$ws.Range("J63").Value = 19371.25
$ws.Range("L63").Value = 19371.25
$ws.Range("N63").Value = -20619.25
$ws.Range("H66").Value = 19371.25
$ws.Range("J66").Value = 19371.25
$ws.Range("L66").Value = 58113.75
$ws.Range("N66").Value = -64353.75
$ws.Range("H96").Value = 4715
$ws.Range("I96").Value = 3667.6667
$ws.Range("K96").Value = 3667.6667
$ws.Range("M96").Value = -2294.6667
$ws.Range("H113").Value = 1491.3529
$ws.Range("J113").Value = 2114.8333
$ws.Range("L113").Value = 6344.499899999999
$ws.Range("N113").Value = -10684.4999
$ws.Range("H122").Value = 4380.696
$ws.Range("I122").Value = 2805.75
$ws.Range("J122").Value = 7980.5713
$ws.Range("K122").Value = 8417.25
$ws.Range("L122").Value = 23941.7139
$ws.Range("M122").Value = -5967.25
$ws.Range("N122").Value = -28841.7139
